# Generate Report for handoff
# Adds two newly-discovered localization files
# (2a4cd2b0-dd70-46fb-ba5d-0b42314b9a92 and 5f27a645-0a48-4464-8bb7-02c82fa3ca64)
# to the Overview/zh-cn/de-de report sheets, refreshes the handoff
# timestamps for the already-tracked files, and keeps the
# ".localization-config" housekeeping row last.

$wb = $excel.ActiveWorkbook

$xlShiftDown = -4121
$xlFormatFromLeftOrAbove = 0

$newHandoffZh = "2016-01-14 04:45:57"
$newHandoffDe = "2016-01-14 04:46:23"

$md1 = "2a4cd2b0-dd70-46fb-ba5d-0b42314b9a92.md"
$md2 = "5f27a645-0a48-4464-8bb7-02c82fa3ca64.md"

$zhXlf1 = "2a4cd2b0-dd70-46fb-ba5d-0b42314b9a92.ea7a97f2a54ebdd002ae622d06eb0c7bb6777733.zh-cn.xlf"
$zhXlf2 = "5f27a645-0a48-4464-8bb7-02c82fa3ca64.4a145be8e693522600a895cc417ee5f2a5e6e80c.zh-cn.xlf"
$deXlf1 = "2a4cd2b0-dd70-46fb-ba5d-0b42314b9a92.ea7a97f2a54ebdd002ae622d06eb0c7bb6777733.de-de.xlf"
$deXlf2 = "5f27a645-0a48-4464-8bb7-02c82fa3ca64.4a145be8e693522600a895cc417ee5f2a5e6e80c.de-de.xlf"

$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/$md1"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/$md2"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/.localization-config"

$zhXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b86db9b4a5937a129b42e4ee38f74529641fbc44/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf1"
$zhXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b86db9b4a5937a129b42e4ee38f74529641fbc44/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf2"
$deXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd03614353adfcdfa94f5fc8e921447858892bc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf1"
$deXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd03614353adfcdfa94f5fc8e921447858892bc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf2"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Make room for the two new rows right before the ".localization-config"
# row (currently row 4), copying formatting from the row above so the
# hyperlink-style font carries over to the new cells.
$wsOverview.Rows.Item(4).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$wsOverview.Rows.Item(5).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)

$wsOverview.Range("A4").Value = $md1
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Range("A5").Value = $md2
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

# Hyperlinks do not shift automatically when rows are inserted, so rebuild
# the whole collection against the final row layout.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $md1Url, [Type]::Missing, [Type]::Missing, "8e32303a-6b08-44be-8c50-6f2f86eba560.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/bb0b23a7-82a4-4308-9794-6ecf61228dbd.md", [Type]::Missing, [Type]::Missing, "bb0b23a7-82a4-4308-9794-6ecf61228dbd.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $md1Url, [Type]::Missing, [Type]::Missing, $md1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $md2Url, [Type]::Missing, [Type]::Missing, $md2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(4).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$wsZh.Rows.Item(5).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)

# Refresh the handoff timestamp for the two already-tracked files.
$wsZh.Range("D2").Value = $newHandoffZh
$wsZh.Range("D3").Value = $newHandoffZh

$wsZh.Range("A4").Value = $md1
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = $zhXlf1
$wsZh.Range("D4").Value = $newHandoffZh
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("A5").Value = $md2
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = $zhXlf2
$wsZh.Range("D5").Value = $newHandoffZh
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/8e32303a-6b08-44be-8c50-6f2f86eba560.md", [Type]::Missing, [Type]::Missing, "8e32303a-6b08-44be-8c50-6f2f86eba560.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b86db9b4a5937a129b42e4ee38f74529641fbc44/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/8e32303a-6b08-44be-8c50-6f2f86eba560.8c127256fd97dbf8ce9d491cd74688131f7f9dc9.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8e32303a-6b08-44be-8c50-6f2f86eba560.8c127256fd97dbf8ce9d491cd74688131f7f9dc9.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/bb0b23a7-82a4-4308-9794-6ecf61228dbd.md", [Type]::Missing, [Type]::Missing, "bb0b23a7-82a4-4308-9794-6ecf61228dbd.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b86db9b4a5937a129b42e4ee38f74529641fbc44/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/bb0b23a7-82a4-4308-9794-6ecf61228dbd.6c0e39badfeb12366594656c2db17e52d01cd394.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "bb0b23a7-82a4-4308-9794-6ecf61228dbd.6c0e39badfeb12366594656c2db17e52d01cd394.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $md1Url, [Type]::Missing, [Type]::Missing, $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhXlf1Url, [Type]::Missing, [Type]::Missing, $zhXlf1)
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $md2Url, [Type]::Missing, [Type]::Missing, $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), $zhXlf2Url, [Type]::Missing, [Type]::Missing, $zhXlf2)
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(4).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$wsDe.Rows.Item(5).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)

$wsDe.Range("D2").Value = $newHandoffDe
$wsDe.Range("D3").Value = $newHandoffDe

$wsDe.Range("A4").Value = $md1
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = $deXlf1
$wsDe.Range("D4").Value = $newHandoffDe
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("A5").Value = $md2
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = $deXlf2
$wsDe.Range("D5").Value = $newHandoffDe
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/8e32303a-6b08-44be-8c50-6f2f86eba560.md", [Type]::Missing, [Type]::Missing, "8e32303a-6b08-44be-8c50-6f2f86eba560.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd03614353adfcdfa94f5fc8e921447858892bc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/8e32303a-6b08-44be-8c50-6f2f86eba560.8c127256fd97dbf8ce9d491cd74688131f7f9dc9.de-de.xlf", [Type]::Missing, [Type]::Missing, "8e32303a-6b08-44be-8c50-6f2f86eba560.8c127256fd97dbf8ce9d491cd74688131f7f9dc9.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/798a8d5db54889e6a3a2b15e486292af6633bec3/e2e/bb0b23a7-82a4-4308-9794-6ecf61228dbd.md", [Type]::Missing, [Type]::Missing, "bb0b23a7-82a4-4308-9794-6ecf61228dbd.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd03614353adfcdfa94f5fc8e921447858892bc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/bb0b23a7-82a4-4308-9794-6ecf61228dbd.6c0e39badfeb12366594656c2db17e52d01cd394.de-de.xlf", [Type]::Missing, [Type]::Missing, "bb0b23a7-82a4-4308-9794-6ecf61228dbd.6c0e39badfeb12366594656c2db17e52d01cd394.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $md1Url, [Type]::Missing, [Type]::Missing, $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deXlf1Url, [Type]::Missing, [Type]::Missing, $deXlf1)
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $md2Url, [Type]::Missing, [Type]::Missing, $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), $deXlf2Url, [Type]::Missing, [Type]::Missing, $deXlf2)
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

Write-Host "Report regenerated for handoff."
